$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + date range) ---
$ws.Range("A8").Value = "Volume 32   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/21/2025  Through  4/27/2025"

# --- Template source cells for style-preserving copies ---
# C14 = style13 shared-string "0" placeholder
# E14 = style13 shared-string "***.*" placeholder
# I14 = style14 numeric template
# K14 = style15 numeric template

# Row 14
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))
$ws.Range("N14").Value = -75
# Row 15
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
# Row 16
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = 16.666666666666
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 41.176470588235
$ws.Range("I16").Value = 72
$ws.Range("J16").Value = 72
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 44
$ws.Range("M16").Value = -17.241379310344
$ws.Range("N16").Value = -74.285714285714
# Row 17
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 34
$ws.Range("H17").Value = -30.612244897959
$ws.Range("I17").Value = 147
$ws.Range("J17").Value = 170
$ws.Range("K17").Value = -13.529411764705
$ws.Range("L17").Value = 5
$ws.Range("M17").Value = 72.941176470588
$ws.Range("N17").Value = -47.122302158273
# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 47
$ws.Range("K18").Value = -31.914893617021
$ws.Range("L18").Value = -43.859649122807
$ws.Range("M18").Value = -8.571428571428
$ws.Range("N18").Value = -89.508196721311
# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -18.75
$ws.Range("I19").Value = 122
$ws.Range("J19").Value = 106
$ws.Range("K19").Value = 15.094339622641
$ws.Range("L19").Value = 14.018691588785
$ws.Range("M19").Value = 35.555555555555
$ws.Range("N19").Value = -0.813008130081
# Row 20
$ws.Range("I14").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 3
$ws.Range("I14").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 1
$ws.Range("K14").Copy($ws.Range("E20"))
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 16
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = -42.857142857142
$ws.Range("L20").Value = -38.461538461538
$ws.Range("M20").Value = 14.285714285714
$ws.Range("N20").Value = -82.978723404255
# Row 21
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 18.518518518518
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -14.912280701754
$ws.Range("I21").Value = 403
$ws.Range("J21").Value = 431
$ws.Range("K21").Value = -6.496519721577
$ws.Range("L21").Value = 2.28426395939
$ws.Range("M21").Value = 25.545171339563
$ws.Range("N21").Value = -64.017857142857
# Row 23
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -40
$ws.Range("F23").Value = 15
$ws.Range("H23").Value = -37.5
$ws.Range("I23").Value = 77
$ws.Range("J23").Value = 89
$ws.Range("K23").Value = -13.483146067415
$ws.Range("L23").Value = 14.925373134328
$ws.Range("M23").Value = 57.142857142857
# Row 24
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 75
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = 4.166666666666
$ws.Range("I24").Value = 285
$ws.Range("J24").Value = 273
$ws.Range("K24").Value = 4.395604395604
$ws.Range("L24").Value = 10.465116279069
$ws.Range("M24").Value = 37.681159420289
# Row 25
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 91.666666666666
$ws.Range("I25").Value = 60
$ws.Range("J25").Value = 59
$ws.Range("K25").Value = 1.694915254237
$ws.Range("L25").Value = -1.639344262295
# Row 26
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 7.692307692307
$ws.Range("F26").Value = 58
$ws.Range("G26").Value = 61
$ws.Range("H26").Value = -4.918032786885
$ws.Range("I26").Value = 183
$ws.Range("J26").Value = 201
$ws.Range("K26").Value = -8.955223880597
$ws.Range("L26").Value = 1.10497237569
$ws.Range("M26").Value = -29.06976744186
# Row 27
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
# Row 28
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("E28").Value = -100
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 18.75
$ws.Range("L28").Value = 0
# Row 29
$ws.Range("I14").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$ws.Range("I14").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K14").Copy($ws.Range("E29"))
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1
$ws.Range("I14").Copy($ws.Range("G29"))
$ws.Range("G29").Value = 1
$ws.Range("K14").Copy($ws.Range("H29"))
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 6
$ws.Range("J29").Value = 8
$ws.Range("K29").Value = -25
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -81.818181818181
# Row 30
$ws.Range("I14").Copy($ws.Range("C30"))
$ws.Range("C30").Value = 1
$ws.Range("I14").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("K14").Copy($ws.Range("E30"))
$ws.Range("E30").Value = 0
$ws.Range("I14").Copy($ws.Range("G30"))
$ws.Range("G30").Value = 1
$ws.Range("K14").Copy($ws.Range("H30"))
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 5
$ws.Range("J30").Value = 5
$ws.Range("L30").Value = -28.571428571428
$ws.Range("M30").Value = -28.571428571428
$ws.Range("N30").Value = -82.758620689655
# Row 31
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 0
# Row 33
$ws.Range("C14").Copy($ws.Range("G33"))
$ws.Range("E14").Copy($ws.Range("H33"))
